$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50 currently holds the (text) phone value "09876543" with a leading
# zero, 0 points. The new data set keeps that original row intact one row
# lower (row 51) and adds a corrected row 50 above it where the phone
# number is stored as a plain number (9876543, no leading zero) with its
# points reset to 0.00.
#
# Insert a fresh row at 50: this pushes the existing row 50 down to 51,
# carrying its original text phone number / blank birthday / 0 points
# along unchanged.
$ws.Rows.Item(50).Insert()

# Populate the new row 50 with the numeric phone number and 0 points.
$ws.Cells.Item(50, 1).Value = 9876543
$ws.Cells.Item(50, 3).Value = 0

# Column B (birthday) is blank for this row. A plain "" assignment clears
# the cell entirely (COM semantics), so use the classic leading-apostrophe
# text marker to commit an actual empty *text* value instead - matching
# the empty-string birthday cells used throughout the rest of the sheet -
# then drop back to the default style so no stray quote-prefix formatting
# is left behind.
$ws.Cells.Item(50, 2).Value = "'"
$ws.Cells.Item(50, 2).Style = "Normal"
